$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("L2").Value = 4.5
$ws.Range("M2").Value = 1.11
$ws.Range("N2").Value = 6.5
$ws.Range("Q2").Value = 1.95
$ws.Range("R2").Value = 1.9

# Row 3
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7

# Row 5
$ws.Range("G5").Value = 1.95
$ws.Range("H5").Value = 3.2
$ws.Range("J5").Value = 2.75
$ws.Range("K5").Value = 1.95
$ws.Range("L5").Value = 5
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 7
$ws.Range("O5").Value = 1.44
$ws.Range("P5").Value = 2.63
$ws.Range("Q5").Value = 1.85
$ws.Range("R5").Value = 2
$ws.Range("S5").Value = 2.4
$ws.Range("T5").Value = 1.53
$ws.Range("U5").Value = 3.7
$ws.Range("V5").Value = 1.27
$ws.Range("W5").Value = 4.5
$ws.Range("X5").Value = 1.18
$ws.Range("Y5").Value = 1.53
$ws.Range("Z5").Value = 2.38
$ws.Range("AA5").Value = 2.2
$ws.Range("AB5").Value = 1.62

# Row 8
$ws.Range("G8").Value = 1.35
$ws.Range("H8").Value = 4.35
$ws.Range("I8").Value = 8.25
$ws.Range("J8").Value = 1.8
$ws.Range("K8").Value = 2.32
$ws.Range("L8").Value = 7.6
$ws.Range("O8").Value = 1.32
$ws.Range("P8").Value = 2.85
$ws.Range("S8").Value = 1.93
$ws.Range("T8").Value = 1.7
$ws.Range("W8").Value = 3.15
$ws.Range("X8").Value = 1.26
$ws.Range("AA8").Value = 2.37
$ws.Range("AB8").Value = 1.45
